$wb = $excel.ActiveWorkbook

# Rename the worksheet from "funneldata" to "SalesFunnelData". Renaming the
# sheet in Excel automatically updates the qualified references in the
# defined name (_xlnm._FilterDatabase) that points at it.
$ws = $wb.Worksheets.Item(1)
$ws.Name = "SalesFunnelData"

# Move the active selection on the sheet from D206 to N205.
$ws.Activate()
$ws.Range("N205").Select()
